# Add a "greedy" column to the GatewaySource sheet, right after the
# reaction_type column (A) and before circular (old B). This shifts the
# existing columns B:H one position to the right (C:I) and sets the new
# B1 header to "greedy".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GatewaySource")

$ws.Columns("B").Insert()
$ws.Range("B1").Value = "greedy"
